$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds a numeric value into a
# text "placeholder" cell (matching the existing style used for A14:N14,
# i.e. right aligned / vertically centered / Andale WT 10pt, no explicit
# number format).  ClearFormats + re-applying alignment & font reproduces
# that exact style so Excel re-uses the existing style record instead of
# creating a new, duplicate one.
# ---------------------------------------------------------------------------
function Set-TextPlaceholder($range, [string]$text) {
    $range.ClearFormats()
    $range.HorizontalAlignment = -4152   # xlRight
    $range.VerticalAlignment = -4108     # xlCenter
    $range.Font.Name = "Andale WT"
    $range.Font.Size = 10
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# Helper: convert a cell that currently holds placeholder text into a
# numeric cell, applying the proper number format used by the other
# numeric cells in the same column group ("#,##0" for raw counts,
# "#,##0.0;""-""#,##0.0" for the percentage-change columns).
# ---------------------------------------------------------------------------
function Set-NumericCount($range, $value) {
    $range.ClearFormats()
    $range.HorizontalAlignment = -4152   # xlRight
    $range.VerticalAlignment = -4108     # xlCenter
    $range.Font.Name = "Andale WT"
    $range.Font.Size = 10
    $range.NumberFormat = "#,##0"
    $range.Value = $value
}

function Set-NumericPct($range, $value) {
    $range.ClearFormats()
    $range.HorizontalAlignment = -4152   # xlRight
    $range.VerticalAlignment = -4108     # xlCenter
    $range.Font.Name = "Andale WT"
    $range.Font.Size = 10
    $range.NumberFormat = '#,##0.0;"-"#,##0.0'
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# Header text updates (report header / volume number)
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "19"
$ws.Range("C9").Characters(46, 8).Text = "5/14/2023"
$ws.Range("C9").Characters(27, 8).Text = "5/8/2023"

# ---------------------------------------------------------------------------
# Row 15 - Murder
# ---------------------------------------------------------------------------
Set-NumericCount $ws.Range("D15") 1
Set-NumericPct   $ws.Range("E15") -100
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -70
$ws.Range("M15").Value = -70
$ws.Range("N15").Value = -90.909090909090

# ---------------------------------------------------------------------------
# Row 16 - Rape
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 90
$ws.Range("I16").Value = 87
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 10.126582278481
$ws.Range("M16").Value = -19.444444444444
$ws.Range("N16").Value = -84.464285714285

# ---------------------------------------------------------------------------
# Row 17 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 45
$ws.Range("I17").Value = 109
$ws.Range("J17").Value = 103
$ws.Range("K17").Value = 5.825242718446
$ws.Range("L17").Value = 18.478260869565
$ws.Range("M17").Value = 57.971014492753
$ws.Range("N17").Value = -64.724919093851

# ---------------------------------------------------------------------------
# Row 18 - Fel. Assault
# ---------------------------------------------------------------------------
Set-NumericCount $ws.Range("C18") 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = 1.639344262295
$ws.Range("L18").Value = 40.909090909090
$ws.Range("M18").Value = -10.144927536231
$ws.Range("N18").Value = -91.621621621621

# ---------------------------------------------------------------------------
# Row 19 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -22.033898305084
$ws.Range("I19").Value = 189
$ws.Range("J19").Value = 232
$ws.Range("K19").Value = -18.534482758620
$ws.Range("L19").Value = -6.435643564356
$ws.Range("M19").Value = 60.169491525423
$ws.Range("N19").Value = -55.424528301886

# ---------------------------------------------------------------------------
# Row 20 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -28.571428571428
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 96
$ws.Range("J20").Value = 134
$ws.Range("K20").Value = -28.358208955223
$ws.Range("L20").Value = 35.211267605633
$ws.Range("M20").Value = 231.034482758621
$ws.Range("N20").Value = -87.723785166240

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold styles, no style transitions required)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 25
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = 12.711864406779
$ws.Range("I21").Value = 546
$ws.Range("J21").Value = 619
$ws.Range("K21").Value = -11.793214862681
$ws.Range("L21").Value = 13.513513513513
$ws.Range("M21").Value = 35.148514851485
$ws.Range("N21").Value = -81.021897810219

# ---------------------------------------------------------------------------
# Row 22 - G.L.A.
# ---------------------------------------------------------------------------
Set-TextPlaceholder $ws.Range("C22") "0"
Set-TextPlaceholder $ws.Range("D22") "0"
Set-TextPlaceholder $ws.Range("E22") "***.*"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300

# ---------------------------------------------------------------------------
# Row 23 - Transit
# ---------------------------------------------------------------------------
Set-NumericCount $ws.Range("C23") 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 14
$ws.Range("K23").Value = -14.285714285714
$ws.Range("L23").Value = 33.333333333333
$ws.Range("M23").Value = 20

# ---------------------------------------------------------------------------
# Row 24 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 27
$ws.Range("E24").Value = 68.75
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 17.808219178082
$ws.Range("I24").Value = 448
$ws.Range("J24").Value = 471
$ws.Range("K24").Value = -4.883227176220
$ws.Range("L24").Value = 62.318840579710
$ws.Range("M24").Value = 109.345794392523

# ---------------------------------------------------------------------------
# Row 25 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 36
$ws.Range("H25").Value = 2.857142857142
$ws.Range("I25").Value = 176
$ws.Range("J25").Value = 166
$ws.Range("K25").Value = 6.024096385542
$ws.Range("L25").Value = 30.370370370370
$ws.Range("M25").Value = -9.743589743589

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
Set-TextPlaceholder $ws.Range("C26") "0"
Set-NumericCount $ws.Range("D26") 3
Set-NumericPct $ws.Range("E26") -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -55.555555555555

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
Set-TextPlaceholder $ws.Range("C27") "0"
$ws.Range("L27").Value = 0

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumericCount $ws.Range("C28") 1
Set-NumericCount $ws.Range("F28") 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = -75
$ws.Range("L28").Value = -77.777777777777
$ws.Range("M28").Value = -60
$ws.Range("N28").Value = -97.619047619047

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-NumericCount $ws.Range("C29") 1
Set-NumericCount $ws.Range("F29") 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = -75
$ws.Range("L29").Value = -75
$ws.Range("M29").Value = -60
$ws.Range("N29").Value = -97.402597402597
